# Add two new columns "I0" and "IF" (columns I and J) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows: row, I-value, J-value
$rows = @(
    @(2, 1, 5),
    @(3, 1, 5),
    @(4, 1, 6),
    @(5, 1, 6),
    @(6, 1, 6),
    @(7, 1, 7),
    @(8, 1, 6),
    @(9, 1, 5),
    @(10, 1, 4),
    @(11, 1, 4),
    @(12, 1, 6),
    @(13, 1, 4),
    @(14, 1, 3),
    @(15, 1, 5),
    @(16, 1, 5),
    @(17, 1, 5),
    @(18, 1, 6),
    @(19, 1, 6),
    @(20, 1, 5),
    @(21, 1, 5),
    @(22, 1, 6),
    @(23, 1, 5),
    @(24, 1, 6),
    @(25, 1, 5),
    @(26, 1, 6),
    @(27, 1, 5),
    @(28, 1, 6),
    @(29, 1, 6),
    @(30, 1, 6),
    @(31, 1, 4),
    @(32, 1, 4),
    @(33, 7, 8),
    @(34, 4, 5),
    @(35, 1, 2)
)

foreach ($entry in $rows) {
    $r = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
